# Applies the "Add_GF_Protein.xlsx" edits:
#  - F2: "angaj2010" -> "Protein_selenium"
#  - C2: protein sequence text gains a leading newline
#  - Row 2 height: 75 -> 90
#  - Column widths for E, H, K, Q change (and bestFit removed); new explicit width for column R
#  - Selection / scroll position moves from F2/B1 to T2/I1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Cell content changes -------------------------------------------------

$seq = "1 marvsanava lvalvsvllt ygccaqspln ytgslakssk aswswlpaka twygaptgag 61 pddnggacgy khtnqypfms mtscgneplf kdgmgcgacy qircvnnkac sgkpetvmit 121 dmnyypvgky hfdlsgtafg amakpgqndk lrhagiidiq  fqrvpcnhpg lnvnfqverg 181 snpnylavlv efanregtvv qmdlmesrng rptgywtamr hswgaiwrmd srrrlqgpfs 241 lrirsesgkt lvakqvipan wrpdtnyrsn vqfr"

$ws.Range("C2").Value = "`n" + $seq
$ws.Range("F2").Value = "Protein_selenium"

# --- Row height -------------------------------------------------------------

$ws.Rows.Item(2).RowHeight = 90

# --- Column widths ------------------------------------------------------
# Stored OOXML widths are ColumnWidth + 5/6, so subtract that offset to hit
# the target stored width as closely as this engine's rounding allows.
$offset = 5/6

$ws.Columns.Item(5).ColumnWidth  = 17.5703125  - $offset   # E
$ws.Columns.Item(8).ColumnWidth  = 21.28515625 - $offset   # H
$ws.Columns.Item(11).ColumnWidth = 22.42578125 - $offset   # K
$ws.Columns.Item(17).ColumnWidth = 18           - $offset  # Q
$ws.Columns.Item(18).ColumnWidth = 25.42578125 - $offset   # R (new explicit width)

# --- View / selection ----------------------------------------------------

$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 9
$win.ScrollRow = 1
$ws.Range("T2").Select()

"done"
